$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 130 <- data from original row 133
$ws.Range("B130").Value2 = 7483247
$ws.Range("C130").Value = 'Ecuador LigaPro Serie A'
$ws.Range("D130").Value2 = 45255.83333333334
$ws.Range("E130").Value = 'Mushuc Runa'
$ws.Range("F130").Value = 'Universidad Catolica del Ecuador'
$ws.Range("G130").Value2 = 0
$ws.Range("H130").Value2 = 2
$ws.Range("I130").Value2 = 0
$ws.Range("J130").Value2 = 1
$ws.Range("K130").Value = 'A'
$ws.Range("L130").Value2 = 3.25
$ws.Range("M130").Value2 = 3.2
$ws.Range("N130").Value2 = 2.25
$ws.Range("O130").Value2 = 3.5
$ws.Range("P130").Value2 = 3.25
$ws.Range("Q130").Value2 = 2.1
$ws.Range("R130").Value2 = 0.5
$ws.Range("S130").Value2 = 1.775
$ws.Range("T130").Value2 = 2.025
$ws.Range("U130").Value2 = 2.5
$ws.Range("V130").Value2 = 1.9
$ws.Range("W130").Value2 = 1.9
$ws.Range("X130").Value2 = -1
$ws.Range("Y130").Value2 = -1
$ws.Range("Z130").Value2 = 1.1
$ws.Range("AA130").Value2 = -1
$ws.Range("AB130").Value2 = 1.025
$ws.Range("AC130").Value2 = -1
$ws.Range("AD130").Value2 = 0.8999999999999999

# Row 131 <- data from original row 132
$ws.Range("B131").Value2 = 7483281
$ws.Range("C131").Value = 'Ecuador LigaPro Serie A'
$ws.Range("D131").Value2 = 45255.83333333334
$ws.Range("E131").Value = 'SD Aucas'
$ws.Range("F131").Value = 'Delfin SC'
$ws.Range("G131").Value2 = 0
$ws.Range("H131").Value2 = 0
$ws.Range("I131").Value2 = 0
$ws.Range("J131").Value2 = 0
$ws.Range("K131").Value = 'D'
$ws.Range("L131").Value2 = 1.909
$ws.Range("M131").Value2 = 3.25
$ws.Range("N131").Value2 = 4.2
$ws.Range("O131").Value2 = 1.909
$ws.Range("P131").Value2 = 3.5
$ws.Range("Q131").Value2 = 4
$ws.Range("R131").Value2 = -0.5
$ws.Range("S131").Value2 = 1.9
$ws.Range("T131").Value2 = 1.9
$ws.Range("U131").Value2 = 2.5
$ws.Range("V131").Value2 = 1.8
$ws.Range("W131").Value2 = 2
$ws.Range("X131").Value2 = -1
$ws.Range("Y131").Value2 = 2.5
$ws.Range("Z131").Value2 = -1
$ws.Range("AA131").Value2 = -1
$ws.Range("AB131").Value2 = 0.8999999999999999
$ws.Range("AC131").Value2 = -1
$ws.Range("AD131").Value2 = 1

# Row 132 <- data from original row 131
$ws.Range("B132").Value2 = 7483081
$ws.Range("C132").Value = 'Ecuador LigaPro Serie A'
$ws.Range("D132").Value2 = 45255.83333333334
$ws.Range("E132").Value = 'Deportivo Cuenca'
$ws.Range("F132").Value = 'El Nacional'
$ws.Range("G132").Value2 = 1
$ws.Range("H132").Value2 = 0
$ws.Range("I132").Value2 = 0
$ws.Range("J132").Value2 = 0
$ws.Range("K132").Value = 'H'
$ws.Range("L132").Value2 = 2.75
$ws.Range("M132").Value2 = 3.25
$ws.Range("N132").Value2 = 2.55
$ws.Range("O132").Value2 = 3
$ws.Range("P132").Value2 = 3.3
$ws.Range("Q132").Value2 = 2.3
$ws.Range("R132").Value2 = 0.25
$ws.Range("S132").Value2 = 1.825
$ws.Range("T132").Value2 = 1.975
$ws.Range("U132").Value2 = 2.75
$ws.Range("V132").Value2 = 2
$ws.Range("W132").Value2 = 1.8
$ws.Range("X132").Value2 = 2
$ws.Range("Y132").Value2 = -1
$ws.Range("Z132").Value2 = -1
$ws.Range("AA132").Value2 = 0.825
$ws.Range("AB132").Value2 = -1
$ws.Range("AC132").Value2 = -1
$ws.Range("AD132").Value2 = 0.8

# Row 133 <- data from original row 130
$ws.Range("B133").Value2 = 7483189
$ws.Range("C133").Value = 'Ecuador LigaPro Serie A'
$ws.Range("D133").Value2 = 45255.83333333334
$ws.Range("E133").Value = 'Independiente del Valle'
$ws.Range("F133").Value = 'Orense'
$ws.Range("G133").Value2 = 2
$ws.Range("H133").Value2 = 2
$ws.Range("I133").Value2 = 1
$ws.Range("J133").Value2 = 0
$ws.Range("K133").Value = 'D'
$ws.Range("L133").Value2 = 1.4
$ws.Range("M133").Value2 = 4.75
$ws.Range("N133").Value2 = 7
$ws.Range("O133").Value2 = 1.4
$ws.Range("P133").Value2 = 4.5
$ws.Range("Q133").Value2 = 8
$ws.Range("R133").Value2 = -1.25
$ws.Range("S133").Value2 = 1.875
$ws.Range("T133").Value2 = 1.925
$ws.Range("U133").Value2 = 2.5
$ws.Range("V133").Value2 = 1.925
$ws.Range("W133").Value2 = 1.875
$ws.Range("X133").Value2 = -1
$ws.Range("Y133").Value2 = 3.5
$ws.Range("Z133").Value2 = -1
$ws.Range("AA133").Value2 = -1
$ws.Range("AB133").Value2 = 0.925
$ws.Range("AC133").Value2 = 0.925
$ws.Range("AD133").Value2 = -1

# Row 142 <- data from original row 144
$ws.Range("B142").Value2 = 7528848
$ws.Range("C142").Value = 'Ecuador LigaPro Serie A'
$ws.Range("D142").Value2 = 45263.83333333334
$ws.Range("E142").Value = 'Emelec'
$ws.Range("F142").Value = 'Deportivo Cuenca'
$ws.Range("G142").Value2 = 2
$ws.Range("H142").Value2 = 1
$ws.Range("I142").Value2 = 0
$ws.Range("J142").Value2 = 1
$ws.Range("K142").Value = 'H'
$ws.Range("L142").Value2 = 1.75
$ws.Range("M142").Value2 = 3.5
$ws.Range("N142").Value2 = 4.2
$ws.Range("O142").Value2 = 2.4
$ws.Range("P142").Value2 = 3.1
$ws.Range("Q142").Value2 = 2.75
$ws.Range("R142").Value2 = -0.25
$ws.Range("S142").Value2 = 2.05
$ws.Range("T142").Value2 = 1.75
$ws.Range("U142").Value2 = 2.25
$ws.Range("V142").Value2 = 1.8
$ws.Range("W142").Value2 = 2
$ws.Range("X142").Value2 = 1.4
$ws.Range("Y142").Value2 = -1
$ws.Range("Z142").Value2 = -1
$ws.Range("AA142").Value2 = 1.05
$ws.Range("AB142").Value2 = -1
$ws.Range("AC142").Value2 = 0.8
$ws.Range("AD142").Value2 = -1

# Row 143 <- data from original row 145
$ws.Range("B143").Value2 = 7528858
$ws.Range("C143").Value = 'Ecuador LigaPro Serie A'
$ws.Range("D143").Value2 = 45263.83333333334
$ws.Range("E143").Value = 'Orense'
$ws.Range("F143").Value = 'SD Aucas'
$ws.Range("G143").Value2 = 1
$ws.Range("H143").Value2 = 2
$ws.Range("I143").Value2 = 1
$ws.Range("J143").Value2 = 1
$ws.Range("K143").Value = 'A'
$ws.Range("L143").Value2 = 2.2
$ws.Range("M143").Value2 = 3.2
$ws.Range("N143").Value2 = 3.2
$ws.Range("O143").Value2 = 1.95
$ws.Range("P143").Value2 = 3.2
$ws.Range("Q143").Value2 = 3.8
$ws.Range("R143").Value2 = -0.5
$ws.Range("S143").Value2 = 1.95
$ws.Range("T143").Value2 = 1.85
$ws.Range("U143").Value2 = 2.25
$ws.Range("V143").Value2 = 1.85
$ws.Range("W143").Value2 = 1.95
$ws.Range("X143").Value2 = -1
$ws.Range("Y143").Value2 = -1
$ws.Range("Z143").Value2 = 2.8
$ws.Range("AA143").Value2 = -1
$ws.Range("AB143").Value2 = 0.8500000000000001
$ws.Range("AC143").Value2 = 0.8500000000000001
$ws.Range("AD143").Value2 = -1

# Row 144 <- data from original row 143
$ws.Range("B144").Value2 = 7528852
$ws.Range("C144").Value = 'Ecuador LigaPro Serie A'
$ws.Range("D144").Value2 = 45263.83333333334
$ws.Range("E144").Value = 'Delfin SC'
$ws.Range("F144").Value = 'Tecnico Universitario'
$ws.Range("G144").Value2 = 2
$ws.Range("H144").Value2 = 2
$ws.Range("I144").Value2 = 1
$ws.Range("J144").Value2 = 0
$ws.Range("K144").Value = 'D'
$ws.Range("L144").Value2 = 2.1
$ws.Range("M144").Value2 = 3.4
$ws.Range("N144").Value2 = 3.1
$ws.Range("O144").Value2 = 2.1
$ws.Range("P144").Value2 = 3.4
$ws.Range("Q144").Value2 = 3.1
$ws.Range("R144").Value2 = -0.25
$ws.Range("S144").Value2 = 1.8
$ws.Range("T144").Value2 = 2
$ws.Range("U144").Value2 = 2.25
$ws.Range("V144").Value2 = 1.9
$ws.Range("W144").Value2 = 1.9
$ws.Range("X144").Value2 = -1
$ws.Range("Y144").Value2 = 2.4
$ws.Range("Z144").Value2 = -1
$ws.Range("AA144").Value2 = -0.5
$ws.Range("AB144").Value2 = 0.5
$ws.Range("AC144").Value2 = 0.8999999999999999
$ws.Range("AD144").Value2 = -1

# Row 145 <- data from original row 142
$ws.Range("B145").Value2 = 7528857
$ws.Range("C145").Value = 'Ecuador LigaPro Serie A'
$ws.Range("D145").Value2 = 45263.83333333334
$ws.Range("E145").Value = 'Universidad Catolica del Ecuador'
$ws.Range("F145").Value = 'Barcelona Guayaquil'
$ws.Range("G145").Value2 = 0
$ws.Range("H145").Value2 = 1
$ws.Range("I145").Value2 = 0
$ws.Range("J145").Value2 = 0
$ws.Range("K145").Value = 'A'
$ws.Range("L145").Value2 = 1.533
$ws.Range("M145").Value2 = 4
$ws.Range("N145").Value2 = 5.5
$ws.Range("O145").Value2 = 1.5
$ws.Range("P145").Value2 = 4.333
$ws.Range("Q145").Value2 = 5.25
$ws.Range("R145").Value2 = -1
$ws.Range("S145").Value2 = 1.8
$ws.Range("T145").Value2 = 2
$ws.Range("U145").Value2 = 3
$ws.Range("V145").Value2 = 1.975
$ws.Range("W145").Value2 = 1.825
$ws.Range("X145").Value2 = -1
$ws.Range("Y145").Value2 = -1
$ws.Range("Z145").Value2 = 4.25
$ws.Range("AA145").Value2 = -1
$ws.Range("AB145").Value2 = 1
$ws.Range("AC145").Value2 = -1
$ws.Range("AD145").Value2 = 0.825

# Row 254 <- data from original row 255
$ws.Range("B254").Value2 = 7773540
$ws.Range("C254").Value = 'Ecuador LigaPro Serie A'
$ws.Range("D254").Value2 = 45437.91666666666
$ws.Range("E254").Value = 'Barcelona Guayaquil'
$ws.Range("F254").Value = 'Tecnico Universitario'
$ws.Range("G254").Value2 = 3
$ws.Range("H254").Value2 = 0
$ws.Range("I254").ClearContents()
$ws.Range("J254").ClearContents()
$ws.Range("K254").Value = 'H'
$ws.Range("L254").Value2 = 1.5
$ws.Range("M254").Value2 = 4
$ws.Range("N254").Value2 = 5.75
$ws.Range("O254").Value2 = 1.285
$ws.Range("P254").Value2 = 4.5
$ws.Range("Q254").Value2 = 12
$ws.Range("R254").Value2 = -1.5
$ws.Range("S254").Value2 = 2.025
$ws.Range("T254").Value2 = 1.775
$ws.Range("U254").Value2 = 2.25
$ws.Range("V254").Value2 = 2
$ws.Range("W254").Value2 = 1.8
$ws.Range("X254").Value2 = 0.2849999999999999
$ws.Range("Y254").Value2 = -1
$ws.Range("Z254").Value2 = -1
$ws.Range("AA254").Value2 = 1.025
$ws.Range("AB254").Value2 = -1
$ws.Range("AC254").Value2 = 1
$ws.Range("AD254").Value2 = -1

# Row 255 <- data from original row 256
$ws.Range("B255").Value2 = 7773537
$ws.Range("C255").Value = 'Ecuador LigaPro Serie A'
$ws.Range("D255").Value2 = 45437.91666666666
$ws.Range("E255").Value = 'LDU Quito'
$ws.Range("F255").Value = 'Universidad Catolica del Ecuador'
$ws.Range("G255").Value2 = 2
$ws.Range("H255").Value2 = 1
$ws.Range("I255").ClearContents()
$ws.Range("J255").ClearContents()
$ws.Range("K255").Value = 'H'
$ws.Range("L255").Value2 = 1.833
$ws.Range("M255").Value2 = 3.5
$ws.Range("N255").Value2 = 3.8
$ws.Range("O255").Value2 = 2.1
$ws.Range("P255").Value2 = 3.4
$ws.Range("Q255").Value2 = 3.1
$ws.Range("R255").Value2 = -0.25
$ws.Range("S255").Value2 = 1.825
$ws.Range("T255").Value2 = 1.975
$ws.Range("U255").Value2 = 2.5
$ws.Range("V255").Value2 = 1.8
$ws.Range("W255").Value2 = 2
$ws.Range("X255").Value2 = 1.1
$ws.Range("Y255").Value2 = -1
$ws.Range("Z255").Value2 = -1
$ws.Range("AA255").Value2 = 0.825
$ws.Range("AB255").Value2 = -1
$ws.Range("AC255").Value2 = 0.8
$ws.Range("AD255").Value2 = -1

# Row 256 <- data from original row 254
$ws.Range("B256").Value2 = 7773538
$ws.Range("C256").Value = 'Ecuador LigaPro Serie A'
$ws.Range("D256").Value2 = 45437.91666666666
$ws.Range("E256").Value = 'Delfin SC'
$ws.Range("F256").Value = 'SD Aucas'
$ws.Range("G256").Value2 = 2
$ws.Range("H256").Value2 = 0
$ws.Range("I256").ClearContents()
$ws.Range("J256").ClearContents()
$ws.Range("K256").Value = 'H'
$ws.Range("L256").Value2 = 3.2
$ws.Range("M256").Value2 = 3.3
$ws.Range("N256").Value2 = 2.1
$ws.Range("O256").Value2 = 3.5
$ws.Range("P256").Value2 = 3.4
$ws.Range("Q256").Value2 = 1.95
$ws.Range("R256").Value2 = 0.5
$ws.Range("S256").Value2 = 1.8
$ws.Range("T256").Value2 = 2
$ws.Range("U256").Value2 = 2.25
$ws.Range("V256").Value2 = 1.775
$ws.Range("W256").Value2 = 2.025
$ws.Range("X256").Value2 = 2.5
$ws.Range("Y256").Value2 = -1
$ws.Range("Z256").Value2 = -1
$ws.Range("AA256").Value2 = 0.8
$ws.Range("AB256").Value2 = -1
$ws.Range("AC256").Value2 = -0.5
$ws.Range("AD256").Value2 = 0.5125
